$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing last data row (177 -> old 178) down by inserting two
# blank rows above the current row 178, so the old row 178 (Clemenuless,
# Primera) becomes row 180 with all of its data/formatting intact.
$ws.Rows("178:179").Insert()

# --- New row 178: Murcott / Especial -------------------------------------
$ws.Cells.Item(178, 1).Value = 11
$ws.Cells.Item(178, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(178, 3).Value = "Bíobío"
$ws.Cells.Item(178, 4).Value = 44911
$ws.Cells.Item(178, 5).Value = 8
$ws.Cells.Item(178, 6).Value = "Fruta"
$ws.Cells.Item(178, 7).Value = 100102
$ws.Cells.Item(178, 8).Value = "Cítricos"
$ws.Cells.Item(178, 9).Value = 100102004
$ws.Cells.Item(178, 10).Value = "Mandarina"
$ws.Cells.Item(178, 11).Value = "Murcott"
$ws.Cells.Item(178, 12).Value = "Especial"
$ws.Cells.Item(178, 13).Value = 50
$ws.Cells.Item(178, 14).Value = 10000
$ws.Cells.Item(178, 15).Value = 10000
$ws.Cells.Item(178, 16).Value = 10000
$ws.Cells.Item(178, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(178, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(178, 19).Value = 556
$ws.Cells.Item(178, 20).Value = 18

# --- New row 179: Murcott / Primera ---------------------------------------
$ws.Cells.Item(179, 1).Value = 11
$ws.Cells.Item(179, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(179, 3).Value = "Bíobío"
$ws.Cells.Item(179, 4).Value = 44911
$ws.Cells.Item(179, 5).Value = 8
$ws.Cells.Item(179, 6).Value = "Fruta"
$ws.Cells.Item(179, 7).Value = 100102
$ws.Cells.Item(179, 8).Value = "Cítricos"
$ws.Cells.Item(179, 9).Value = 100102004
$ws.Cells.Item(179, 10).Value = "Mandarina"
$ws.Cells.Item(179, 11).Value = "Murcott"
$ws.Cells.Item(179, 12).Value = "Primera"
$ws.Cells.Item(179, 13).Value = 100
$ws.Cells.Item(179, 14).Value = 9000
$ws.Cells.Item(179, 15).Value = 9000
$ws.Cells.Item(179, 16).Value = 9000
$ws.Cells.Item(179, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(179, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(179, 19).Value = 500
$ws.Cells.Item(179, 20).Value = 18
